$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "36.522.81"
$ws.Cells.Item(2, 5).Value = "  -0.73%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.100.00"
$ws.Cells.Item(3, 5).Value = "  +9.77%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "253.23"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.13%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.656"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -6.74%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "47.58"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +2.23%  "

# Row 9
$ws.Cells.Item(9, 2).Value = "Cardano"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.381"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.35%  "

# Row 10
$ws.Cells.Item(10, 2).Value = "OKB"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "59.97"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +2.43%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.89%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +0.26%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "14.53"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.15%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.404.60"

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.837"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.72%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.099.93"
$ws.Cells.Item(16, 5).Value = "  +9.66%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -1.07%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "36.512.87"
$ws.Cells.Item(18, 5).Value = "  -0.71%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "72.73"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -2.96%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0831"
$ws.Cells.Item(20, 5).Value = "  -3.73%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.23"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.34%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "240.27"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -4.35%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.22"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.84%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -0.06%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.51"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -4.59%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "171.30"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.88%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +13.75%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +3.95%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.99"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -9.49%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "28.87"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +60.91%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -5.29%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -1.54%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0615"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.35%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "LidoDAOToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.44"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +21.53%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.986"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +12.94%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0916"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.97%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.08%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.86"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.09%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "4.08"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -5.80%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -11.05%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "VeChain"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0223"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.73%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "ARBITRUM"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.17"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +6.24%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "98.07"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -6.70%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -5.13%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "15.95"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -9.57%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "1.332.24"
$ws.Cells.Item(46, 5).Value = "  -1.37%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +3.62%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.11"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +9.93%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "2.304.86"
$ws.Cells.Item(49, 5).Value = "  +10.11%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +1.81%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.25"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -5.29%  "
